# Update scraped "想去人数" (want-to-go count) figures in the 杭州-漫展信息
# workbook across the 展览, 演出 and 全部类型 sheets, matching the freshly
# scraped data snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsAll  = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet ---
$wsExpo.Range("F6").Value  = 2896
$wsExpo.Range("F7").Value  = 1691
$wsExpo.Range("F8").Value  = 1922
$wsExpo.Range("F11").Value = 782
$wsExpo.Range("F12").Value = 929
$wsExpo.Range("F13").Value = 186
$wsExpo.Range("F14").Value = 394
$wsExpo.Range("F19").Value = 6963
$wsExpo.Range("F20").Value = 262
$wsExpo.Range("F21").Value = 1697
$wsExpo.Range("F22").Value = 174
$wsExpo.Range("F25").Value = 347
$wsExpo.Range("F26").Value = 279
$wsExpo.Range("F27").Value = 72
$wsExpo.Range("F28").Value = 1115
$wsExpo.Range("F29").Value = 928
$wsExpo.Range("F31").Value = 104
$wsExpo.Range("F33").Value = 801
$wsExpo.Range("F35").Value = 165
$wsExpo.Range("F37").Value = 148
$wsExpo.Range("F41").Value = 245
$wsExpo.Range("F43").Value = 181

# --- 演出 sheet ---
$wsShow.Range("F2").Value = 17
$wsShow.Range("F6").Value = 7

# --- 全部类型 sheet ---
$wsAll.Range("F5").Value  = 17
$wsAll.Range("F9").Value  = 2896
$wsAll.Range("F10").Value = 1691
$wsAll.Range("F11").Value = 1922
$wsAll.Range("F14").Value = 782
$wsAll.Range("F16").Value = 929
$wsAll.Range("F17").Value = 186
$wsAll.Range("F18").Value = 394
$wsAll.Range("F22").Value = 6963
$wsAll.Range("F23").Value = 262
$wsAll.Range("F24").Value = 1697
$wsAll.Range("F25").Value = 7
$wsAll.Range("F26").Value = 174
$wsAll.Range("F29").Value = 347
$wsAll.Range("F30").Value = 279
$wsAll.Range("F31").Value = 72
$wsAll.Range("F32").Value = 1115
$wsAll.Range("F33").Value = 0
$wsAll.Range("F35").Value = 104
$wsAll.Range("F36").Value = 801
$wsAll.Range("F38").Value = 165
$wsAll.Range("F40").Value = 148
$wsAll.Range("F44").Value = 245
$wsAll.Range("F49").Value = 181
